$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped on Tue Sep  5 23:59:34 UTC 2023

$ws.Range("D2").Value = '25.880.74'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '1.638.50'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.10'
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2589'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06434'
$ws.Range("E9").Value = '  +1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.38'
$ws.Range("E10").Value = '  +4.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07801'
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.269'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.642.26'
$ws.Range("E13").Value = '  +0.53%  '
$ws.Range("D14").Value = '1.865.86'
$ws.Range("E14").Value = '  +0.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5601'
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("D16").Value = '0.0₅7672'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.32'
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("D18").Value = '25.896.82'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.11'
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.388'
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.951'
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.157'
$ws.Range("E23").Value = '  +2.17%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.793'
$ws.Range("E25").Value = '  -4.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '138.16'
$ws.Range("E26").Value = '  -2.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1233'
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.839'
$ws.Range("E28").Value = '  +1.55%  '
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04974'
$ws.Range("E31").Value = '  +2.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.301'
$ws.Range("E32").Value = '  +0.97%  '
$ws.Range("E33").Value = '  +2.18%  '
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("E35").Value = '  +1.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9049'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("E37").Value = '  +1.38%  '
$ws.Range("E38").Value = '  +0.20%  '
$ws.Range("D39").Value = '1.137.02'
$ws.Range("E39").Value = '  +1.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01574'
$ws.Range("E40").Value = '  +1.17%  '
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.71'
$ws.Range("E42").Value = '  +2.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.480'
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8029'
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("E45").Value = '  -2.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.65'
$ws.Range("E46").Value = '  +1.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4251'
$ws.Range("E47").Value = '  -3.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.787'
$ws.Range("E48").Value = '  +2.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05062'
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("E51").Value = '  +0.41%  '
